# Bianca_Project / OverviewSheet.xlsx edit
# - Remove Bianca Saia's e-mail address from F12 (content now cleared, style kept)
#   Removing the only use of the "bianca@divanbleu.com" shared string causes it
#   to be dropped from sharedStrings.xml on save, shifting "Rent" (used by B1)
#   down one shared-string index automatically.
# - Move the sheet's active cell / selection from E15 to G8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear Bianca Saia's e-mail cell (row 12, column F) while keeping its format.
$ws.Range("F12").ClearContents()

# Update the saved selection shown when the workbook is reopened.
$ws.Range("G8").Select()
